# ---------------------------------------------------------------------------
# Applies the "Functional Requirements" updates described by the commit:
#   1. "Allow admins and authorised staff to change timetables" gains
#      " and rotas", followed by four new sub-bullets about bus lines/routes.
#   2. "Staff assigned to particular routes" is split into two runs (the
#      second part, "particular routes", was flagged by the grammar checker).
#   3. The "search routes" bullet gets a clarifying bracketed phrase and
#      "correspond to" becomes "go through".
#
# NOTE on technique: all *text* mutations are performed first. Only once the
# final wording is in place do we touch (toggle-and-restore Bold on) the
# sub-ranges that the diff shows as separate <w:r> runs -- this persuades the
# engine's run-coalescing pass to keep those runs distinct without altering
# any visible formatting (the toggle nets out to a no-op).
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- Change 1a: "...change timetables" -> "...change timetables and rotas"
$d.Content.Find.Execute(
    "Allow admins and authorised staff to change timetables",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Allow admins and authorised staff to change timetables and rotas", 2) | Out-Null

# --- Change 1b: insert the four new ilvl=1 bullets right after that bullet
# (NOTE: Paragraph.Index is unreliable in this engine -- it returns
# inconsistent/off-by-one values -- so the paragraph's 1-based position is
# tracked by walking $d.Paragraphs ourselves instead of trusting .Index.)
$timetablesParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd() -eq "Allow admins and authorised staff to change timetables and rotas") {
        $timetablesParaIndex = $i
        break
    }
}

$newBullets = @(
    "Create bus lines and routes",
    "Allocate staff to lines",
    "Make sure the staff are working the hours they were contracted for",
    "Make sure no driver is assigned to two routes at the same time"
)
$d.Paragraphs.Item($timetablesParaIndex).Range.InsertAfter("`r" + ($newBullets -join "`r"))

# New paragraphs land directly after the "...and rotas" bullet; promote them
# to the second list level (ilvl=1) to match their siblings.
for ($i = 0; $i -lt $newBullets.Count; $i++) {
    $bulletPara = $d.Paragraphs.Item($timetablesParaIndex + 1 + $i)
    $bulletPara.Range.ListFormat.ListLevelNumber = 2
}

# --- Change 3: reword the "search routes" bullet
$d.Content.Find.Execute(
    "search for terms and find the routes that correspond to",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "search for terms [town names and bus stop names] and find the routes that go through", 2) | Out-Null

# ---------------------------------------------------------------------------
# All wording is final now -- recreate the run boundaries the diff shows.
# ---------------------------------------------------------------------------
function Split-Run($rangeToIsolate) {
    $rangeToIsolate.Bold = 1
    $rangeToIsolate.Bold = 0
}

# Change 1: " and rotas" as its own run
$r1 = $d.Content
$r1.Find.Execute(" and rotas", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Split-Run $r1

# Change 2: "particular routes" as its own run
$r2 = $d.Content
$r2.Find.Execute("particular routes", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Split-Run $r2

# Change 3: " [town names and bus stop names]" and "go through" as their own runs
$r3 = $d.Content
$r3.Find.Execute(" [town names and bus stop names]", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Split-Run $r3

$r4 = $d.Content
$r4.Find.Execute("go through", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Split-Run $r4

Write-Output "Done."
